$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells keep their original text representation
# (values such as "1.012" or "337.51" would otherwise be auto-converted
# to numbers by Excel, losing the trailing zero / text formatting), so
# the Price column cells are forced to Text format before assignment.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.639.76'
$ws.Range("E2").Value = '  +0.55%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.117.72'
$ws.Range("E3").Value = '  +1.17%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.012'
$ws.Range("E4").Value = '  +1.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '337.51'
$ws.Range("E5").Value = '  +2.19%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.011'
$ws.Range("E6").Value = '  +0.94%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5248'
$ws.Range("E7").Value = '  +0.55%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4549'
$ws.Range("E8").Value = '  +2.56%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '54.53'
$ws.Range("E9").Value = '  +1.09%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.09104'
$ws.Range("E10").Value = '  +1.79%  '

$ws.Range("E11").Value = '  +1.96%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '24.47'
$ws.Range("E12").Value = '  +0.72%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.118.48'
$ws.Range("E13").Value = '  +1.29%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.856'
$ws.Range("E14").Value = '  +2.50%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.134'
$ws.Range("E15").Value = '  +5.48%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001173'
$ws.Range("E16").Value = '  +4.57%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '97.18'
$ws.Range("E17").Value = '  +1.31%  '

$ws.Range("E18").Value = '  +0.92%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06686'
$ws.Range("E19").Value = '  +1.16%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.44'
$ws.Range("E20").Value = '  +1.75%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.010'
$ws.Range("E21").Value = '  +0.84%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.307'
$ws.Range("E22").Value = '  +0.63%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '30.703.86'
$ws.Range("E23").Value = '  +0.62%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.92'
$ws.Range("E24").Value = '  +4.98%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.359'
$ws.Range("E25").Value = '  +1.90%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.361.85'
$ws.Range("E26").Value = '  +1.23%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.36'
$ws.Range("E27").Value = '  +0.55%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '164.88'
$ws.Range("E28").Value = '  +0.82%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.554'
$ws.Range("E29").Value = '  -0.46%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '134.79'
$ws.Range("E30").Value = '  +2.47%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.209'
$ws.Range("E31").Value = '  +1.96%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1074'
$ws.Range("E32").Value = '  +0.41%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.644'
$ws.Range("E33").Value = '  -0.55%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.353'
$ws.Range("E34").Value = '  +3.30%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.949'
$ws.Range("E35").Value = '  +1.18%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '10.68'
$ws.Range("E36").Value = '  +5.45%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.869'
$ws.Range("E37").Value = '  +7.42%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02638'
$ws.Range("E38").Value = '  +3.21%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06866'
$ws.Range("E39").Value = '  +1.05%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2327'
$ws.Range("E40").Value = '  +3.23%  '

$ws.Range("E41").Value = '  -0.40%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6898'
$ws.Range("E42").Value = '  +0.22%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.261'
$ws.Range("E43").Value = '  +0.91%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.89'
$ws.Range("E44").Value = '  +6.55%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6485'
$ws.Range("E45").Value = '  +2.57%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.311'
$ws.Range("E46").Value = '  +5.08%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000370'
$ws.Range("E47").Value = '  +22.63%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.699'
$ws.Range("E48").Value = '  +2.05%  '

$ws.Range("E49").Value = '  +0.85%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '83.36'
$ws.Range("E50").Value = '  +2.15%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.194'
$ws.Range("E51").Value = '  -3.60%  '
